# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect the newer scrape captured in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# "展览" sheet — rows keyed directly by their row number.
$exhibitionUpdates = @(
    @{ Row = 3;  New = 5408 },
    @{ Row = 7;  New = 619 },
    @{ Row = 9;  New = 1059 },
    @{ Row = 11; New = 1488 },
    @{ Row = 12; New = 4432 },
    @{ Row = 17; New = 3519 },
    @{ Row = 25; New = 45 }
)

# "全部类型" sheet — same events, shifted down by one row vs. "展览".
$allTypesUpdates = @(
    @{ Row = 4;  New = 5408 },
    @{ Row = 8;  New = 619 },
    @{ Row = 10; New = 1059 },
    @{ Row = 12; New = 1488 },
    @{ Row = 13; New = 4432 },
    @{ Row = 18; New = 3519 },
    @{ Row = 26; New = 45 }
)

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($u in $exhibitionUpdates) {
    $wsExhibition.Cells.Item($u.Row, 6).Value = $u.New
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($u in $allTypesUpdates) {
    $wsAllTypes.Cells.Item($u.Row, 6).Value = $u.New
}
